# austin_away_passing.xlsx — "multiple cleaning changes, added sql schema"
#
# The diff shows a new "Match ID" column inserted as the new column A of the
# single worksheet: every existing column shifts one place to the right
# (A->B, B->C, ... AC->AD), the merged header cells (H1:L1 etc.) shift with
# them, and the new column A is filled with the header "Match ID" (row 3,
# bold, no border) and the literal value 1 for every data row (rows 4-20;
# bold for the visible rows 4-19, unstyled for the hidden summary row 20).
# The view's selection also moves from A3:AC3 to the single cell G30.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new column before A — shifts all existing data/styles/merged
# cells one column to the right without touching their contents.
$ws.Columns("A:A").Insert()

# New header cell for the inserted column (row 3 holds the column headers).
$ws.Range("A3").Value = "Match ID"
$ws.Range("A3").Font.Bold = $true

# Fill every data row (including the hidden totals row 20) with the match id.
$ws.Range("A4:A20").Value = 1

# Visible player rows (4-19) get the same bold styling as the rest of row 3's
# header-ish "Player ID" column; the hidden summary row (20) stays unstyled.
$ws.Range("A4:A19").Font.Bold = $true

# Writing into row 20 (hidden) can otherwise stamp a stray custom row height;
# re-autofit it so the row keeps its original (default) height metadata.
$ws.Rows("20:20").EntireRow.AutoFit()

# Match the saved view state: selection on G30 (top-left scroll position is
# not exposed through this host's object model).
$ws.Range("G30").Select()
